# 10.11.2024 - wykonanie mini pomocy oraz szablonow plikow do importu
#
# Updates the "Sheet1" task list:
#   - D2: 0 -> 75
#   - D3: 0 -> 100
#   - D24: 0 -> 100
#   - E2: new note "dodać ramki do zestawień jeszcze. Postarać się sformatować całość"
#   - Column E gets a custom width
#   - Selection moves to G8
# F1/G1/H1 are formulas and recalc automatically from the D-column edits.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 75
$ws.Range("D3").Value = 100
$ws.Range("D24").Value = 100

$ws.Range("E2").Value = "dodać ramki do zestawień jeszcze. Postarać się sformatować całość"

# Widen column E (closest reachable width to the authored 13.5703125 chars
# given this host's pixel-grid rounding of ColumnWidth).
$ws.Columns.Item(5).ColumnWidth = 12.6666666666667

$ws.Range("G8").Select() | Out-Null
